# Refresh cryptocurrency Price (col D) and Volume(1h) (col E) figures.
# Column D holds plain-text numbers (t="inlineStr" in the sheet XML); a
# bare $ws.Range(...).Value = "183.85" would let Excel's COM layer auto-
# coerce that into a real Double. Prefixing with an apostrophe (exactly like
# typing '183.85 into a General-formatted cell in the Excel UI) keeps it text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.091.47"
$ws.Range("E2").Value = "  -4.80%  "
$ws.Range("D3").Value = "3.314.07"
$ws.Range("E3").Value = "  -5.88%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'183.85"
$ws.Range("E5").Value = "  -9.06%  "
$ws.Range("D6").Value = "'527.82"
$ws.Range("E6").Value = "  -4.38%  "
$ws.Range("D7").Value = "'0.605"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.315.26"
$ws.Range("E8").Value = "  -5.61%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'0.626"
$ws.Range("E10").Value = "  -5.03%  "
$ws.Range("D11").Value = "'60.32"
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("E12").Value = "  -6.62%  "
$ws.Range("D13").Value = "'0.0000262"
$ws.Range("E13").Value = "  -4.81%  "
$ws.Range("D14").Value = "'9.21"
$ws.Range("E14").Value = "  -6.47%  "
$ws.Range("D15").Value = "3.821.80"
$ws.Range("E15").Value = "  -6.46%  "
$ws.Range("D16").Value = "'0.119"
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("D17").Value = "3.302.09"
$ws.Range("E17").Value = "  -5.91%  "
$ws.Range("D18").Value = "'17.81"
$ws.Range("E18").Value = "  -4.50%  "
$ws.Range("D19").Value = "63.906.39"
$ws.Range("E19").Value = "  -4.78%  "
$ws.Range("E20").Value = "  -6.59%  "
$ws.Range("D21").Value = "'0.965"
$ws.Range("E21").Value = "  -6.92%  "
$ws.Range("D22").Value = "'374.35"
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("D23").Value = "'11.36"
$ws.Range("E23").Value = "  -5.72%  "
$ws.Range("D24").Value = "'3.78"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("D25").Value = "'81.11"
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("E26").Value = "  +4.84%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'2.70"
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("D29").Value = "'11.62"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").Value = "'8.49"
$ws.Range("E30").Value = "  -4.96%  "
$ws.Range("D31").Value = "'6.97"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").Value = "'29.02"
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("D33").Value = "'650.26"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").Value = "'11.43"
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").Value = "'0.107"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("D36").Value = "'59.34"
$ws.Range("E36").Value = "  -7.26%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'37.22"
$ws.Range("E39").Value = "  -6.29%  "
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "2.940.45"
$ws.Range("E42").Value = "  -5.36%  "
$ws.Range("D43").Value = "0.0₃0675"
$ws.Range("E43").Value = "  -4.34%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  -10.43%  "
$ws.Range("D46").Value = "'2.96"
$ws.Range("E46").Value = "  +12.74%  "
$ws.Range("D47").Value = "'0.0402"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "'2.65"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").Value = "'2.63"
$ws.Range("E49").Value = "  -6.47%  "
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'2.98"
$ws.Range("E51").Value = "  +1.15%  "
